$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '72.746.01'
$ws.Range("E2").Value = '  +4.38%  '

$ws.Range("D3").Value = '3.965.60'
$ws.Range("E3").Value = '  +1.61%  '

$ws.Range("E4").Value = '  +0.23%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '587.89'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +10.90%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '158.26'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +8.50%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.683'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.04%  '

$ws.Range("E8").Value = '  -0.15%  '

$ws.Range("E9").Value = '  +3.14%  '

$ws.Range("E10").Value = '  +3.21%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '54.22'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.10%  '

$ws.Range("E12").Value = '  +2.33%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '10.86'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +4.63%  '

$ws.Range("D14").Value = '4.606.48'
$ws.Range("E14").Value = '  +2.61%  '

$ws.Range("D15").Value = '3.974.27'
$ws.Range("E15").Value = '  +2.54%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.26'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +9.84%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.05'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.74%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '20.32'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.94%  '

$ws.Range("E19").Value = '  +0.57%  '

$ws.Range("D20").Value = '72.475.38'
$ws.Range("E20").Value = '  +4.15%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '434.12'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.85%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.68'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +14.40%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '95.88'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.15%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.43'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.45%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '14.33'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.17%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '4.38'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +21.71%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.12'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.00%  '

$ws.Range("E28").Value = '  +1.44%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.94'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.99%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '36.35'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.64%  '

$ws.Range("E31").Value = '  +5.77%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '13.65'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.98%  '

$ws.Range("E33").Value = '  +4.62%  '

$ws.Range("B34").Value = 'Bittensor'
$ws.Range("C34").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '681.04'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.59%  '

$ws.Range("B35").Value = 'InjectiveProtocol'
$ws.Range("C35").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '48.51'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.60%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '69.63'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +9.85%  '

$ws.Range("B37").Value = 'PEPE'
$ws.Range("C37").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D37").Value = '0.0₃0872'
$ws.Range("E37").Value = '  +9.77%  '

$ws.Range("B38").Value = 'TheGraph'
$ws.Range("C38").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.435'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.71%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.38'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.24%  '

$ws.Range("E40").Value = '  -0.46%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.997'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.38%  '

$ws.Range("E42").Value = '  +4.23%  '

$ws.Range("E43").Value = '  +0.61%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '10.84'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +14.62%  '

$ws.Range("E45").Value = '  +2.89%  '

$ws.Range("E46").Value = '  +2.31%  '

$ws.Range("E47").Value = '  -0.60%  '

$ws.Range("E48").Value = '  +3.25%  '

$ws.Range("E49").Value = '  +3.12%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.39'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +6.53%  '

$ws.Range("D51").Value = '2.790.22'
$ws.Range("E51").Value = '  +11.77%  '
